$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.254.71'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '1.894.56'
$ws.Range('E3').Value = '  +1.56%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.661'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.86%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.83%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.349'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.81'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +12.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0715'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.86%  '
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('D13').Value = '2.171.31'
$ws.Range('E13').Value = '  +1.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.05'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.67%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.694'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.14%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.902.35'
$ws.Range('E16').Value = '  +1.99%  '
$ws.Range('E17').Value = '  +1.87%  '
$ws.Range('D18').Value = '35.259.85'
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '71.97'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.69%  '
$ws.Range('E20').Value = '  +2.43%  '
$ws.Range('E21').Value = '  -0.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.41'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.13%  '
$ws.Range('E23').Value = '  +1.41%  '
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('E25').Value = '  +1.33%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.32'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +22.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '170.06'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.39'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.37%  '
$ws.Range('E29').Value = '  +2.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.126'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.36%  '
$ws.Range('E31').Value = '  +2.13%  '
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('E33').Value = '  -0.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.929'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +13.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.07'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.61%  '
$ws.Range('E36').Value = '  -4.70%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.02'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.12%  '
$ws.Range('E38').Value = '  +1.06%  '
$ws.Range('E39').Value = '  -1.47%  '
$ws.Range('E40').Value = '  +2.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '15.97'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0629'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.33%  '
$ws.Range('E43').Value = '  -1.39%  '
$ws.Range('D44').Value = '1.335.67'
$ws.Range('E44').Value = '  -0.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.36'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '47.83'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +37.14%  '
$ws.Range('E47').Value = '  -0.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.77'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.48'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.73'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.60%  '
$ws.Range('D51').Value = '2.082.28'
$ws.Range('E51').Value = '  +1.50%  '
